# Lambda.pptx - "class slides and fixing missing video"
#
# Slide 2 (sldId 300), shape "Content Placeholder 2": remove the bullet
# paragraph "The starting code errors if you load it, but don't worry..."
# while leaving the rest of the bulleted list (and the trailing blank
# paragraph) untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$sh = $s.Shapes.Item("Content Placeholder 2")
$tr = $sh.TextFrame.TextRange

# Locate the paragraph that needs to be removed by matching its text,
# rather than assuming a fixed index, so the script is resilient to
# minor paragraph-count differences.
$paraCount = $tr.Paragraphs().Count
for ($i = $paraCount; $i -ge 1; $i--) {
    $para = $tr.Paragraphs($i, 1)
    if ($para.Text -like "The starting code errors if you load it*") {
        $para.Delete()
        break
    }
}
